$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ShallowVerification")
$ws2 = $wb.Worksheets.Item("SpaceEx")

# --- Sheet1 (ShallowVerification): header I1 + data rows 2-23 ---
$ws1.Cells.Item(1, 9).Value = "Result"
$ws1.Cells.Item(2, 1).Value = 1
$ws1.Cells.Item(2, 4).Value = 1
$ws1.Cells.Item(2, 5).Value = 2176
$ws1.Cells.Item(2, 6).Value = 0.028
$ws1.Cells.Item(2, 7).Value = 0.023
$ws1.Cells.Item(2, 8).Value = 0
$ws1.Cells.Item(2, 9).Value = "UNSATISFIABLE"
$ws1.Cells.Item(3, 1).Value = 2
$ws1.Cells.Item(3, 4).Value = 1
$ws1.Cells.Item(3, 5).Value = 2214
$ws1.Cells.Item(3, 6).Value = 0.024
$ws1.Cells.Item(3, 7).Value = 0.028
$ws1.Cells.Item(3, 8).Value = 0
$ws1.Cells.Item(3, 9).Value = "SATISFIABLE"
$ws1.Cells.Item(4, 1).Value = 3
$ws1.Cells.Item(4, 4).Value = 1
$ws1.Cells.Item(4, 5).Value = 11173
$ws1.Cells.Item(4, 6).Value = 0.086
$ws1.Cells.Item(4, 7).Value = 0.112
$ws1.Cells.Item(4, 8).Value = 1
$ws1.Cells.Item(4, 9).Value = "UNSATISFIABLE"
$ws1.Cells.Item(5, 1).Value = 4
$ws1.Cells.Item(5, 4).Value = 1
$ws1.Cells.Item(5, 5).Value = 11593
$ws1.Cells.Item(5, 6).Value = 0.085
$ws1.Cells.Item(5, 7).Value = 0.283
$ws1.Cells.Item(5, 8).Value = 1
$ws1.Cells.Item(5, 9).Value = "SATISFIABLE"
$ws1.Cells.Item(6, 1).Value = 5
$ws1.Cells.Item(6, 4).Value = 1
$ws1.Cells.Item(6, 5).Value = 1737
$ws1.Cells.Item(6, 6).Value = 0.038
$ws1.Cells.Item(6, 7).Value = 0.016
$ws1.Cells.Item(6, 8).Value = 1
$ws1.Cells.Item(6, 9).Value = "UNSATISFIABLE"
$ws1.Cells.Item(7, 1).Value = 6
$ws1.Cells.Item(7, 4).Value = 1
$ws1.Cells.Item(7, 5).Value = 1730
$ws1.Cells.Item(7, 6).Value = 0.045
$ws1.Cells.Item(7, 7).Value = 0.026
$ws1.Cells.Item(7, 8).Value = 1
$ws1.Cells.Item(7, 9).Value = "UNSATISFIABLE"
$ws1.Cells.Item(8, 1).Value = 7
$ws1.Cells.Item(8, 4).Value = 1
$ws1.Cells.Item(8, 5).Value = 1980
$ws1.Cells.Item(8, 6).Value = 0.045
$ws1.Cells.Item(8, 7).Value = 0.036
$ws1.Cells.Item(8, 8).Value = 2
$ws1.Cells.Item(8, 9).Value = "UNSATISFIABLE"
$ws1.Cells.Item(9, 1).Value = 8
$ws1.Cells.Item(9, 4).Value = 1
$ws1.Cells.Item(9, 5).Value = 1980
$ws1.Cells.Item(9, 6).Value = 0.044
$ws1.Cells.Item(9, 7).Value = 0.027
$ws1.Cells.Item(9, 8).Value = 2
$ws1.Cells.Item(9, 9).Value = "UNSATISFIABLE"
$ws1.Cells.Item(10, 1).Value = 9
$ws1.Cells.Item(10, 4).Value = 1
$ws1.Cells.Item(10, 5).Value = 2695
$ws1.Cells.Item(10, 6).Value = 0.058
$ws1.Cells.Item(10, 7).Value = 0.014
$ws1.Cells.Item(10, 8).Value = 3
$ws1.Cells.Item(10, 9).Value = "UNSATISFIABLE"
$ws1.Cells.Item(11, 1).Value = 10
$ws1.Cells.Item(11, 4).Value = 1
$ws1.Cells.Item(11, 5).Value = 2695
$ws1.Cells.Item(11, 6).Value = 0.025
$ws1.Cells.Item(11, 7).Value = 0.013
$ws1.Cells.Item(11, 8).Value = 3
$ws1.Cells.Item(11, 9).Value = "UNSATISFIABLE"
$ws1.Cells.Item(12, 1).Value = 11
$ws1.Cells.Item(12, 4).Value = 1
$ws1.Cells.Item(12, 5).Value = 1355
$ws1.Cells.Item(12, 6).Value = 0.014
$ws1.Cells.Item(12, 7).Value = 0.016
$ws1.Cells.Item(12, 8).Value = 0
$ws1.Cells.Item(12, 9).Value = "UNSATISFIABLE"
$ws1.Cells.Item(13, 1).Value = 12
$ws1.Cells.Item(13, 4).Value = 1
$ws1.Cells.Item(13, 5).Value = 1355
$ws1.Cells.Item(13, 6).Value = 0.018
$ws1.Cells.Item(13, 7).Value = 0.019
$ws1.Cells.Item(13, 8).Value = 0
$ws1.Cells.Item(13, 9).Value = "UNSATISFIABLE"
$ws1.Cells.Item(14, 1).Value = 13
$ws1.Cells.Item(14, 4).Value = 1
$ws1.Cells.Item(14, 5).Value = 526
$ws1.Cells.Item(14, 6).Value = 0.017
$ws1.Cells.Item(14, 7).Value = 0.018
$ws1.Cells.Item(14, 8).Value = 0
$ws1.Cells.Item(14, 9).Value = "UNSATISFIABLE"
$ws1.Cells.Item(15, 1).Value = 14
$ws1.Cells.Item(15, 4).Value = 1
$ws1.Cells.Item(15, 5).Value = 519
$ws1.Cells.Item(15, 6).Value = 0.012
$ws1.Cells.Item(15, 7).Value = 0.015
$ws1.Cells.Item(15, 8).Value = 0
$ws1.Cells.Item(15, 9).Value = "UNSATISFIABLE"
$ws1.Cells.Item(16, 1).Value = 15
$ws1.Cells.Item(16, 4).Value = 1
$ws1.Cells.Item(16, 5).Value = 1250
$ws1.Cells.Item(16, 6).Value = 0.019
$ws1.Cells.Item(16, 7).Value = 0.016
$ws1.Cells.Item(16, 8).Value = 0
$ws1.Cells.Item(16, 9).Value = "UNSATISFIABLE"
$ws1.Cells.Item(17, 1).Value = 16
$ws1.Cells.Item(17, 4).Value = 1
$ws1.Cells.Item(17, 5).Value = 1235
$ws1.Cells.Item(17, 6).Value = 0.02
$ws1.Cells.Item(17, 7).Value = 0.013
$ws1.Cells.Item(17, 8).Value = 0
$ws1.Cells.Item(17, 9).Value = "UNSATISFIABLE"
$ws1.Cells.Item(18, 1).Value = 17
$ws1.Cells.Item(18, 4).Value = 1
$ws1.Cells.Item(18, 5).Value = 1796
$ws1.Cells.Item(18, 6).Value = 0.02
$ws1.Cells.Item(18, 7).Value = 0.01
$ws1.Cells.Item(18, 8).Value = 0
$ws1.Cells.Item(18, 9).Value = "UNSATISFIABLE"
$ws1.Cells.Item(19, 1).Value = 18
$ws1.Cells.Item(19, 4).Value = 1
$ws1.Cells.Item(19, 5).Value = 1796
$ws1.Cells.Item(19, 6).Value = 0.018
$ws1.Cells.Item(19, 7).Value = 0.013
$ws1.Cells.Item(19, 8).Value = 0
$ws1.Cells.Item(19, 9).Value = "UNSATISFIABLE"
$ws1.Cells.Item(20, 1).Value = 19
$ws1.Cells.Item(20, 4).Value = 1
$ws1.Cells.Item(20, 5).Value = 9418
$ws1.Cells.Item(20, 6).Value = 0.059
$ws1.Cells.Item(20, 7).Value = 0.147
$ws1.Cells.Item(20, 8).Value = 2
$ws1.Cells.Item(20, 9).Value = "UNSATISFIABLE"
$ws1.Cells.Item(21, 1).Value = 20
$ws1.Cells.Item(21, 4).Value = 1
$ws1.Cells.Item(21, 5).Value = 9485
$ws1.Cells.Item(21, 6).Value = 0.071
$ws1.Cells.Item(21, 7).Value = 0.125
$ws1.Cells.Item(21, 8).Value = 2
$ws1.Cells.Item(21, 9).Value = "SATISFIABLE"
$ws1.Cells.Item(22, 1).Value = 21
$ws1.Cells.Item(22, 4).Value = 1
$ws1.Cells.Item(22, 5).Value = 2350
$ws1.Cells.Item(22, 6).Value = 0.025
$ws1.Cells.Item(22, 7).Value = 0.012
$ws1.Cells.Item(22, 8).Value = 1
$ws1.Cells.Item(22, 9).Value = "UNSATISFIABLE"
$ws1.Cells.Item(23, 1).Value = 22
$ws1.Cells.Item(23, 4).Value = 1
$ws1.Cells.Item(23, 5).Value = 1403
$ws1.Cells.Item(23, 6).Value = 0.018
$ws1.Cells.Item(23, 7).Value = 0.012
$ws1.Cells.Item(23, 8).Value = 1
$ws1.Cells.Item(23, 9).Value = "UNSATISFIABLE"

# --- Sheet2 (SpaceEx): new data for rows 20-21 ---
$ws2.Cells.Item(20, 4).Value = 7.77
$ws2.Cells.Item(20, 5).Value = 37
$ws2.Cells.Item(20, 6).Value = 6
$ws2.Cells.Item(20, 7).Value = 2
$ws2.Cells.Item(20, 8).Value = 4
$ws2.Cells.Item(20, 10).Value = 667
$ws2.Cells.Item(21, 4).Value = 8.24
$ws2.Cells.Item(21, 5).Value = 38
$ws2.Cells.Item(21, 6).Value = 6
$ws2.Cells.Item(21, 7).Value = 2
$ws2.Cells.Item(21, 8).Value = 4
$ws2.Cells.Item(21, 10).Value = 693

# --- Column width for new column I on sheet1 ---
$ws1.Columns.Item(9).ColumnWidth = 13.5

# --- Selection / active sheet changes ---
$ws2.Range("D22").Select()
$ws1.Activate()
$ws1.Range("K18").Select()
